$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Feuil1"

$ws.Range("A1").Value = "nom"
$ws.Range("B1").Value = "prenom"
$ws.Range("C1").Value = "age"
$ws.Range("D1").Value = "sexe"

$ws.Range("E1").Select()
